$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update package metadata values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/outcome"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Elements": move the ele-1/ext-1 constraint text ---
# from the "Extension" row (row 2) to the "Extension.extension" row (row 4)
$elements = $wb.Worksheets.Item("Elements")
$constraintText = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Range("AI2").Value = ""
$elements.Range("AI4").Value = $constraintText
